# "Visualizacion de celular y doble copia"
# - Adjusts the datetime display of B153 on "Proveedores" so it shows the
#   time portion too (matches the style used by the newly-added row).
# - Adds a new transaction row (154) on "Proveedores" for "Prueba 3" /
#   Ahorro / Abono of 50000, mirroring the existing "doble copia" entry.
# - Refreshes the "Resumen" rollup sheet: every provider's Id bumps from 13
#   to 14, and "Prueba 3"'s Total Ahorro / Saldo pick up the new 50000
#   abono (50000 -> 100000, -50000 -> -100000).

$wb = $excel.ActiveWorkbook

# ---- Sheet "Proveedores" ----------------------------------------------
$ws1 = $wb.Worksheets.Item("Proveedores")

# B153 keeps its value, only the display format changes to include time.
$ws1.Range("B153").NumberFormat = "yyyy-mm-dd h:mm:ss"

# New row 154.
$ws1.Range("A154").Value = 153
$ws1.Range("B154").Value = 45901
$ws1.Range("B154").NumberFormat = "yyyy-mm-dd"
$ws1.Range("C154").Value = "Prueba 3"
$ws1.Range("D154").Value = "Ahorro"
$ws1.Range("E154").Value = "Abono"
$ws1.Range("F154").Value = 50000

# ---- Sheet "Resumen" ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("Resumen")

foreach ($r in 2..9) {
    $ws2.Range("A$r").Value = 14
}

$ws2.Range("D9").Value = 100000
$ws2.Range("E9").Value = -100000
